$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) After the first paragraph ("Play Catch of the Day Reeling 'Em In
#    free online", Heading1) insert a new paragraph:
#      <empty run><bold "Meta description"><": Read our review ...">
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Style = "Normal"

$p2 = $d.Paragraphs.Item(2)
$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:r/>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
           "<w:r><w:t>: Read our review of Catch of the Day Reeling 'Em In online slot game and play for free. Features, payouts, and RTP explained.</w:t></w:r>" +
           "</w:p>"
$insRange = $d.Range($p2.Range.Start, $p2.Range.End)
$insRange.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Remove the bold "Play Catch of the Day Reeling 'Em In free
#    online" paragraph near the end of the document. (The Heading1
#    paragraph at the very top of the doc has the same wording, so
#    match on the "Normal"-styled, non-heading copy specifically.)
# ------------------------------------------------------------------
$titleText = "Play Catch of the Day Reeling 'Em In free online"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $titleText -and $p.Range.ParagraphStyle.NameLocal -eq "Normal") {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ------------------------------------------------------------------
# 3) Replace the text of the remaining (italic) paragraph with the
#    new image-prompt copy, keeping its italic formatting intact.
#    Range.Text is used (instead of Find/Replace) so that Word's
#    smart-quote autocorrect does not mangle the straight quotes.
# ------------------------------------------------------------------
$oldText = "Read our review of Catch of the Day Reeling 'Em In online slot game and play for free. Features, payouts, and RTP explained."
$newText = 'Create a cartoon-style feature image for the game "Catch of the Day Reeling ‘Em In" featuring a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior holding a fishing rod and wearing a fishing hat, vest, and boots. He should be standing on a boat in the middle of the ocean, with fish jumping up in the background. The image should convey a sense of fun, excitement, and adventure.'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $oldText) {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Text = $newText
        break
    }
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
